# Correct emission factors in local_grid_parameters.xlsx
# The "emission_factor" column (E) was stored in g CO2eq/kWh by mistake;
# convert every value to kg CO2eq/kWh by dividing by 1000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$corrections = @{
    3  = 0.25
    4  = 0.207
    5  = 0.425
    6  = 0.445
    7  = 0.189
    8  = 0.338
    9  = 0.9
    10 = 0.353
    11 = 0.662
    12 = 0.207
    13 = 0.054
    14 = 0.135
    15 = 0.248
    16 = 0.664
    17 = 0.138
    18 = 0.065
    19 = 0.069
    20 = 0.251
    21 = 0.356
    22 = 0.441
    23 = 0.102
    24 = 0.789
    25 = 0.31
    26 = 0.293
    27 = 0.248
    28 = 0.137
    29 = 0.111
    30 = 0.013
    31 = 0.25
    34 = 0.019
}

foreach ($row in $corrections.Keys) {
    $ws.Range("E$row").Value = $corrections[$row]
}

# Cosmetic view-state touch-ups matching the saved file (selection cursor
# moved to F1 before the sheet was re-saved).
$ws.Range("F1").Select()
